$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.691.97"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "1.694.49"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3951"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4058"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08837"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.055"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001318"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "1.699.28"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07025"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.015"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "24.679.57"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.265"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.368"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "136.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.186"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.585"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("D32").Value = "1.884.68"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.067"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08598"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.151"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2741"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09245"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02723"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.465"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7637"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7168"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.583"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.213"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.326"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07977"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.10%  "
